$wb = $excel.ActiveWorkbook

# Overview sheet: update the row for de573e3a-be4a-435b-918e-27473b76cc5b.md (row 3)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-12-14 04:12:15"

# zh-cn sheet: update Status / Latest Handoff Datetime for de573e3a row (row 3)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "2016-03-14 04:12:13"

# de-de sheet: update Status / Latest Handoff Datetime for de573e3a row (row 3)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "2016-03-14 04:12:15"
